$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '56.667.06'
$c.Style = 'Normal'
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.40%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '3.005.63'
$c.Style = 'Normal'
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '513.54'
$c.Style = 'Normal'
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.99%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '137.11'
$c.Style = 'Normal'
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.36%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.429'
$c.Style = 'Normal'
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.92%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '7.26'
$c.Style = 'Normal'
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.10%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = 'Normal'
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.12%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.364'
$c.Style = 'Normal'
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.15%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '3.482.87'
$c.Style = 'Normal'
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.10%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.71%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '26.41'
$c.Style = 'Normal'
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.0000166'
$c.Style = 'Normal'
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.85%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '56.563.31'
$c.Style = 'Normal'
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.56%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '6.12'
$c.Style = 'Normal'
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.11%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '2.996.40'
$c.Style = 'Normal'
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.40%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '13.07'
$c.Style = 'Normal'
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = '  +2.17%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '8.00'
$c.Style = 'Normal'
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.74%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '323.36'
$c.Style = 'Normal'
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.29%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '0.499'
$c.Style = 'Normal'
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.37%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '64.04'
$c.Style = 'Normal'
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '3.089.82'
$c.Style = 'Normal'
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.29%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '0.163'
$c.Style = 'Normal'
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.66%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0880'
$c.Style = 'Normal'
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.93%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '6.55'
$c.Style = 'Normal'
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.47%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '7.18'
$c.Style = 'Normal'
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.45%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.80'
$c.Style = 'Normal'
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.47%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '1.20'
$c.Style = 'Normal'
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.55%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '20.51'
$c.Style = 'Normal'
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.30%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '152.52'
$c.Style = 'Normal'
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.95%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '4.53'
$c.Style = 'Normal'
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.91%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '5.78'
$c.Style = 'Normal'
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.34%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '25.33'
$c.Style = 'Normal'
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.93%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '1.23'
$c.Style = 'Normal'
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.87%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.0663'
$c.Style = 'Normal'
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.65%  '
$c.Style = 'Normal'

$ws.Cells.Item(40, 2).Value = 'OKB'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '36.60'
$c.Style = 'Normal'
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.95%  '
$c.Style = 'Normal'

$ws.Cells.Item(41, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.06%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '3.82'
$c.Style = 'Normal'
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.39%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.650'
$c.Style = 'Normal'
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.28%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.53%  '
$c.Style = 'Normal'

$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '2.162.81'
$c.Style = 'Normal'
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = '@'
$c.Value = '  -4.91%  '
$c.Style = 'Normal'

$ws.Cells.Item(46, 2).Value = 'Cosmos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '6.10'
$c.Style = 'Normal'
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.70%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '0.942'
$c.Style = 'Normal'
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = '@'
$c.Value = '  -4.09%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '0.0241'
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.51%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '19.42'
$c.Style = 'Normal'
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.16%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '1.75'
$c.Style = 'Normal'
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.07%  '
$c.Style = 'Normal'

$ws.Cells.Item(51, 2).Value = 'TheGraph'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '0.181'
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.86%  '
$c.Style = 'Normal'
